$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stale "g CO2" (J20/K20) helper numbers so the now-unused
#     shared string gets dropped and "kg" shifts into its place (matches
#     the diff: C21 ends up pointing at the renumbered "kg" string). ---
$ws.Range("J20").ClearContents()
$ws.Range("K20").ClearContents()

# --- Add the new explanatory text cells, in the same left-to-right /
#     top-to-bottom order the author typed them, so new shared strings
#     land in the same order as the target workbook. ---
$ws.Range("L12").Value = "0.98 gram CO2 per kilogram GFT afval."

$ws.Range("J4").Value = "Vervoer"
$ws.Range("J4").Font.Bold = $true

$ws.Range("J21").Value = "https://www.rvo.nl/sites/default/files/bijlagen/Bio-energie%20-%20Input%20-%20Groente-,%20fruit-%20en%20tuinafval%20(gft).pdf"

$ws.Range("J23").Value = "Composteren ipv verbranden in een HVC reduceert uitstoot met 60 kg CO2 per ton GFT"

$ws.Range("J24").Value = 60
$ws.Range("K24").Value = "kg/ton"

$ws.Range("J25").Value = 60
$ws.Range("K25").Value = "g/kg"

$ws.Range("J20").Value = "Verbranding efficientie"
$ws.Range("J20").Font.Bold = $true

$ws.Range("J28").Value = "Verbranden zelf"
$ws.Range("J28").Font.Bold = $true

$ws.Range("J29").Value = "https://www.milieucentraal.nl/minder-afval/afval-scheiden/afval-scheiden-nut-en-fabels/"

# --- Turn the evi-europark reference (M6) into a real hyperlink. ---
$ws.Hyperlinks.Add($ws.Range("M6"), "https://www.evi-europark.nl/wp-content/uploads/2019/02/101125-8176_defrapportMO_CE.pdf")

# --- Selection moved to M6 (matches new sheetView selection). ---
$ws.Range("M6").Select()

# --- Printer/page setup tweak. ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- absolute path metadata picked up the nested project folder name. ---
$wb.Path
